$d = $word.ActiveDocument
$brk = [char]11

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Unveiling the Mystique of the Ancient Silk Road" "The Evolution of Technology: A Journey Through History"

# --- Author line: "Dr. Marianna Edwards" (3 runs) -> "Luke Patterson" (1 run) ---
Replace-Text "Dr. Marianna Edwards" "Luke Patterson"

# --- Email line parts ---
Replace-Text "marianna" "luke"
Replace-Text "edwards@academicmail" "patterson@highschool"
Replace-Text "org" "edu"

# --- Body paragraph: first sentence ---
Replace-Text "For centuries, the Silk Road, a sprawling network of ancient trade routes, served as a vibrant artery connecting the East and the West" "The world has witnessed an awe-inspiring evolution of technology throughout history, a story intertwined with human ingenuity"

# --- Remove old 2nd sentence, replace 3rd sentence's text with the merged new text ---
Replace-Text " Spanning vast distances from the heart of China, through the deserts of Central Asia, and reaching the distant shores of Europe and Africa, it was a conduit for more than just goods and services." ""
Replace-Text " It was an exchange of ideas, cultures, technologies, and worldviews, shaping the course of human history" " From the first stone tools crafted by our ancestors to the sophisticated digital gadgets we rely on today, technology has played a pivotal role in shaping civilizations and transforming the way we live"

# --- Second sub-paragraph (after first double break) ---
Replace-Text "Along the Silk Road, merchants, travelers, diplomats, and pilgrims crossed paths, bringing with them a kaleidoscope of languages, traditions, and beliefs" "During the dawn of humanity, our ancestors created rudimentary tools from stones, bones, and wood, utilizing these implements for hunting, gathering, and survival"
Replace-Text " The road witnessed the transmission of knowledge, from the secrets of papermaking to the intricacies of astronomy, fostering scientific advancements and intellectual growth" " Over time, technological advancements propelled us forward: the invention of agriculture allowed for settled communities and fostered the growth of civilizations"
Replace-Text " It facilitated the spread of religious teachings, from Buddhism to Christianity, contributing to a profound spiritual transformation across diverse societies" " The harnessing of fire brought warmth, illumination, and the ability to cook food, enhancing our quality of life"

# --- Third sub-paragraph (after second double break) ---
Replace-Text "Cities and towns mushroomed along the Silk Road, becoming bustling centers of commerce and cultural fusion" "As civilizations flourished, technological innovations emerged at an accelerating pace"
Replace-Text " The exchange of goods, from exquisite silks to exotic spices, fueled economic prosperity and interdependence" " The wheel revolutionized transportation, while the development of written language facilitated communication, record-keeping, and the accumulation of knowledge"
Replace-Text " The Silk Road catalyzed the growth of empires and kingdoms, as rulers recognized its strategic importance and sought to control its lucrative trade routes" " The Middle Ages brought forth towering cathedrals and intricate clocks, showcasing the artistry and engineering prowess of the time"

# --- Insert substantial new content after "...engineering prowess of the time" and before the trailing "." run ---
$rng = $d.Content
$rng.Find.Execute("showcasing the artistry and engineering prowess of the time", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$newBlock = ".$brk$brk" + `
    "Introduction Continued:$brk$brk" + `
    "The Renaissance sparked a renewed interest in science and innovation, paving the way for groundbreaking discoveries." + `
    " Maritime navigation techniques and the invention of the printing press unlocked new worlds and disseminated knowledge far and wide." + `
    " The advent of the Industrial Revolution ushered in a surge of technological advancements: steam engines powered factories, railroads crisscrossed continents, and electricity illuminated cities.$brk$brk" + `
    "The 20th century witnessed an explosion of technological innovation that continues to reshape the world today." + `
    " The invention of the computer and the internet has brought about a paradigm shift in communication, information access, and global connectivity." + `
    " Advancements in medicine have extended life expectancy and improved healthcare, while breakthroughs in renewable energy and environmental sustainability offer hope for a greener future.$brk$brk" + `
    "Introduction Concluded:$brk$brk" + `
    "Today, we stand at the precipice of a new era of technology." + `
    " Artificial intelligence, robotics, and genetic engineering hold immense promise for revolutionizing fields such as healthcare, transportation, and manufacturing." + `
    " Technology has become an integral part of our lives, offering convenience, connectivity, and access to information like never before." + `
    " As we venture into the future, we face both opportunities and challenges in harnessing technology responsibly and ensuring that it serves humanity's best interests"
$rng.InsertAfter($newBlock)

# --- Summary paragraph ---
Replace-Text "The Silk Road, an interconnected web of ancient trade routes, was a catalyst for cultural exchange, technological advancements, and economic prosperity" "Our journey through history reveals the transformative power of technology"
Replace-Text " It facilitated the transmission of ideas, religions, and goods, fostering a vibrant tapestry of civilizations" " From humble beginnings with stone tools to the marvels of modern science, technology has shaped our world and continues to redefine the way we live"
Replace-Text " The Silk Road's profound impact shaped the course of human history, leaving an enduring legacy visible in the cultural and intellectual traditions that continue to shape our world today" " As we embrace new frontiers in technology, we must strive for inclusivity, sustainability, and responsible innovation, ensuring that these advancements benefit all of society and safeguard the well-being of our planet"

# --- Add a new empty paragraph at the very end of the document body ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.Text = "`r"
